# Add a new "Effort" worksheet with historical fishing-effort estimates,
# place it after the existing sheets, and make it the active tab.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Effort"

$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "n_days_tot"
$ws.Range("A2").Value = 2001
$ws.Range("B2").Value = 3388

$ws.Range("E31").Select() | Out-Null
